# Auto-generated edit applying weekly Fruta/Hortaliza price update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44391
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15500
$ws.Range("P2").Value = 861
$ws.Range("D3").Value = 44742
$ws.Range("I3").Value = 'Segunda'
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 16000
$ws.Range("M3").Value = 15500
$ws.Range("P3").Value = 861
$ws.Range("D4").Value = 44377
$ws.Range("J4").Value = 100
$ws.Range("M4").Value = 17600
$ws.Range("P4").Value = 978
$ws.Range("D5").Value = 44748
$ws.Range("J5").Value = 250
$ws.Range("D6").Value = 44755
$ws.Range("J6").Value = 160
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 17500
$ws.Range("P6").Value = 972
$ws.Range("D9").Value = 44433
$ws.Range("K9").Value = 17000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 17500
$ws.Range("P9").Value = 972
$ws.Range("D10").Value = 44433
$ws.Range("H10").Value = 'Cultivar IV Región'
$ws.Range("I10").Value = 'Tercera'
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 14000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 14500
$ws.Range("N10").Value = '$/bandeja 18 kilos'
$ws.Range("O10").Value = 'Provincia de Limarí'
$ws.Range("P10").Value = 806
$ws.Range("Q10").Value = 18
$ws.Range("D11").Value = 44435
$ws.Range("H11").Value = 'Cultivar IV Región'
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 17500
$ws.Range("N11").Value = '$/bandeja 18 kilos'
$ws.Range("O11").Value = 'Provincia de Limarí'
$ws.Range("P11").Value = 972
$ws.Range("Q11").Value = 18
$ws.Range("D12").Value = 44435
$ws.Range("H12").Value = 'Cultivar IV Región'
$ws.Range("I12").Value = 'Tercera'
$ws.Range("J12").Value = 120
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = 14500
$ws.Range("N12").Value = '$/bandeja 18 kilos'
$ws.Range("O12").Value = 'Provincia de Limarí'
$ws.Range("P12").Value = 806
$ws.Range("Q12").Value = 18
$ws.Range("D13").Value = 44405
$ws.Range("I13").Value = 'Segunda'
$ws.Range("J13").Value = 140
$ws.Range("D15").Value = 44526
$ws.Range("H15").Value = 'Cultivar XV región'
$ws.Range("I15").Value = 'Primera'
$ws.Range("K15").Value = 5000
$ws.Range("L15").Value = 5500
$ws.Range("M15").Value = 5250
$ws.Range("N15").Value = '$/caja 10 kilos'
$ws.Range("O15").Value = 'Región de Arica y Parinacota'
$ws.Range("P15").Value = 525
$ws.Range("Q15").Value = 10
$ws.Range("D16").Value = 44526
$ws.Range("H16").Value = 'Cultivar XV región'
$ws.Range("I16").Value = 'Segunda'
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 4000
$ws.Range("L16").Value = 4500
$ws.Range("M16").Value = 4250
$ws.Range("N16").Value = '$/caja 10 kilos'
$ws.Range("O16").Value = 'Región de Arica y Parinacota'
$ws.Range("P16").Value = 425
$ws.Range("Q16").Value = 10
$ws.Range("D17").Value = 44526
$ws.Range("I17").Value = 'Tercera'
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 3500
$ws.Range("M17").Value = 3250
$ws.Range("P17").Value = 325
$ws.Range("D18").Value = 44454
$ws.Range("K18").Value = 19000
$ws.Range("L18").Value = 20000
$ws.Range("M18").Value = 19500
$ws.Range("P18").Value = 1083
$ws.Range("D19").Value = 44554
$ws.Range("J19").Value = 200
$ws.Range("L19").Value = 6000
$ws.Range("M19").Value = 5500
$ws.Range("P19").Value = 550
$ws.Range("D20").Value = 44769
$ws.Range("H20").Value = 'Cultivar IV Región'
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 140
$ws.Range("K20").Value = 17000
$ws.Range("L20").Value = 18000
$ws.Range("M20").Value = 17500
$ws.Range("N20").Value = '$/bandeja 18 kilos'
$ws.Range("O20").Value = 'Provincia de Limarí'
$ws.Range("P20").Value = 972
$ws.Range("Q20").Value = 18
$ws.Range("D21").Value = 44783
$ws.Range("H21").Value = 'Cultivar IV Región'
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 150
$ws.Range("K21").Value = 17000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 17500
$ws.Range("N21").Value = '$/bandeja 18 kilos'
$ws.Range("O21").Value = 'Provincia de Limarí'
$ws.Range("P21").Value = 972
$ws.Range("Q21").Value = 18
$ws.Range("D22").Value = 44221
$ws.Range("H22").Value = 'Cultivar XV región'
$ws.Range("J22").Value = 140
$ws.Range("K22").Value = 5000
$ws.Range("L22").Value = 6000
$ws.Range("M22").Value = 5500
$ws.Range("N22").Value = '$/caja 10 kilos'
$ws.Range("O22").Value = 'Región de Arica y Parinacota'
$ws.Range("P22").Value = 550
$ws.Range("Q22").Value = 10
$ws.Range("D23").Value = 44363
$ws.Range("H23").Value = 'Cultivar IV Región'
$ws.Range("J23").Value = 140
$ws.Range("K23").Value = 14000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 14500
$ws.Range("N23").Value = '$/bandeja 18 kilos'
$ws.Range("O23").Value = 'Provincia de Limarí'
$ws.Range("P23").Value = 806
$ws.Range("Q23").Value = 18
$ws.Range("D24").Value = 44211
$ws.Range("H24").Value = 'Cultivar XV región'
$ws.Range("J24").Value = 140
$ws.Range("K24").Value = 4500
$ws.Range("L24").Value = 5000
$ws.Range("M24").Value = 4750
$ws.Range("N24").Value = '$/caja 10 kilos'
$ws.Range("O24").Value = 'Región de Arica y Parinacota'
$ws.Range("P24").Value = 475
$ws.Range("Q24").Value = 10
$ws.Range("D25").Value = 44776
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 17000
$ws.Range("L25").Value = 18000
$ws.Range("M25").Value = 17500
$ws.Range("P25").Value = 972
$ws.Range("D26").Value = 44533
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 6000
$ws.Range("L26").Value = 7000
$ws.Range("M26").Value = 6500
$ws.Range("P26").Value = 650
$ws.Range("D27").Value = 44533
$ws.Range("H27").Value = 'Cultivar XV región'
$ws.Range("I27").Value = 'Segunda'
$ws.Range("J27").Value = 120
$ws.Range("K27").Value = 4000
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = 4500
$ws.Range("N27").Value = '$/caja 10 kilos'
$ws.Range("O27").Value = 'Región de Arica y Parinacota'
$ws.Range("P27").Value = 450
$ws.Range("Q27").Value = 10
$ws.Range("D28").Value = 44762
$ws.Range("K28").Value = 15000
$ws.Range("L28").Value = 16000
$ws.Range("M28").Value = 15500
$ws.Range("P28").Value = 861
$ws.Range("D29").Value = 44757
$ws.Range("H29").Value = 'Cultivar XV región'
$ws.Range("I29").Value = 'Primera'
$ws.Range("J29").Value = 150
$ws.Range("K29").Value = 6000
$ws.Range("L29").Value = 6500
$ws.Range("M29").Value = 6250
$ws.Range("N29").Value = '$/caja 10 kilos'
$ws.Range("O29").Value = 'Región de Arica y Parinacota'
$ws.Range("P29").Value = 625
$ws.Range("Q29").Value = 10
$ws.Range("D30").Value = 44412
